# actualizacion base de datos
# se agrega archivo venv para ingresar credenciales, y rutas configurables, correos y trm

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the product list in column A.
# Row order changes: "zapatos" (row2) is replaced by "mouse", "gorra" (row5) is
# replaced by "cartera", and a new row "alfombra" is appended.
$ws.Range("A1").Value = "items"
$ws.Range("A2").Value = "mouse"
$ws.Range("A3").Value = "xbox"
$ws.Range("A4").Value = "play station"
$ws.Range("A5").Value = "cartera"
$ws.Range("A6").Value = "alfombra"

# Leave the selection where the author left it when saving.
$ws.Range("D13").Select()
